$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("A1").Value = "Horas trabajadas"
$ws.Range("B1").Value = "Personas residentes viviendas familiares"
$ws.Range("C1").Value = "Comarca nombre"
$ws.Range("D1").Value = "Comarca código"
$ws.Range("E1").Value = "Provincia código"
$ws.Range("F1").Value = "Aragón"
$ws.Range("G1").Value = "Provincia nombre"

# Row 2
$ws.Range("A2").Value = "iaest-measure:horas-trabajadas"
$ws.Range("B2").Value = "iaest-measure:personas-residentes-viviendas-familiares"
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("D2").Value = "null"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "sdmx-dimension:refArea"

# Row 3
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "dim"
$ws.Range("G3").Value = "dim"

# Row 4
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "URI-comarca"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("G4").Value = "URI-Provincia"
